$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.991.18'
$ws.Range("E2").Value = '  -0.04%  '

$ws.Range("D3").Value = '1.918.10'
$ws.Range("E3").Value = '  +0.92%  '

$ws.Range("E4").Value = '  -0.73%  '

$ws.Range("D5").Value = '''325.10'
$ws.Range("E5").Value = '  +0.21%  '

$ws.Range("E6").Value = '  -0.54%  '

$ws.Range("D7").Value = '''0.4600'
$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").Value = '''0.3821'
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '''0.07729'
$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("D10").Value = '''0.9806'
$ws.Range("E10").Value = '  +1.63%  '

$ws.Range("D11").Value = '''22.26'
$ws.Range("E11").Value = '  +1.61%  '

$ws.Range("D12").Value = '1.891.31'
$ws.Range("E12").Value = '  -0.31%  '

$ws.Range("D13").Value = '''5.690'
$ws.Range("E13").Value = '  -0.14%  '

$ws.Range("D14").Value = '''6.962'
$ws.Range("E14").Value = '  -0.43%  '

$ws.Range("D15").Value = '''0.07034'
$ws.Range("E15").Value = '  -0.32%  '

$ws.Range("E16").Value = '  -0.44%  '

$ws.Range("D17").Value = '''84.14'
$ws.Range("E17").Value = '  +0.35%  '

$ws.Range("D18").Value = '''0.000009495'
$ws.Range("E18").Value = '  -0.45%  '

$ws.Range("E19").Value = '  -1.37%  '

$ws.Range("D20").Value = '''1.001'
$ws.Range("E20").Value = '  -0.33%  '

$ws.Range("D21").Value = '28.976.83'
$ws.Range("E21").Value = '  -0.12%  '

$ws.Range("D22").Value = '''5.338'
$ws.Range("E22").Value = '  -1.27%  '

$ws.Range("D23").Value = '''10.95'
$ws.Range("E23").Value = '  +0.08%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = '''2.082'
$ws.Range("E24").Value = '  -0.23%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '''158.22'
$ws.Range("E25").Value = '  +0.61%  '

$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = '''19.03'
$ws.Range("E26").Value = '  -0.47%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = '''5.689'
$ws.Range("E27").Value = '  +0.71%  '

$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").Value = '''118.02'
$ws.Range("E28").Value = '  +0.71%  '

$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '''1.855'
$ws.Range("E29").Value = '  +2.32%  '

$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '''0.09323'
$ws.Range("E30").Value = '  +0.33%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '''0.8677'
$ws.Range("E31").Value = '  +1.95%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''5.111'
$ws.Range("E32").Value = '  +0.59%  '

$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").Value = '''1.252'
$ws.Range("E33").Value = '  -0.38%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '''3.053'
$ws.Range("E34").Value = '  -0.81%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.05714'
$ws.Range("E35").Value = '  +0.76%  '

$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").Value = '''1.151'
$ws.Range("E36").Value = '  -0.81%  '

$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").Value = '''1.001'
$ws.Range("E37").Value = '  -0.60%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.02043'
$ws.Range("E38").Value = '  -0.15%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '''3.055'
$ws.Range("E39").Value = '  +12.44%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '''7.515'
$ws.Range("E40").Value = '  +0.69%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '''0.5511'
$ws.Range("E41").Value = '  -0.61%  '

$ws.Range("B42").Value = 'PEPE'
$ws.Range("C42").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D42").Value = '''0.000003019'
$ws.Range("E42").Value = '  +5.64%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '''0.1753'
$ws.Range("E43").Value = '  -0.03%  '

$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '''9.384'
$ws.Range("E44").Value = '  +2.10%  '

$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = '''2.244'
$ws.Range("E45").Value = '  +9.22%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.5179'
$ws.Range("E46").Value = '  -0.25%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''11.25'
$ws.Range("E47").Value = '  +0.00%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '''0.06899'
$ws.Range("E48").Value = '  +1.60%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '''1.780'
$ws.Range("E49").Value = '  -0.16%  '

$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '''110.47'
$ws.Range("E50").Value = '  +0.08%  '

$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").Value = '''1.002'
$ws.Range("E51").Value = '  -0.29%  '
